$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("A1").Value = "特采编号"
$ws.Range("B1").Value = "特采名称"
$ws.Range("C1").Value = "地图"
$ws.Range("D1").Value = "兑换奖励"
$ws.Range("E1").Value = "描述"
$ws.Range("D2").Value = "5000新币"

# Used range extends through column I, row 5 - fill in the rest of the
# grid so every cell in A1:I5 is materialised (even if empty) with the
# centered style.
$ws.Range("A1:I5").HorizontalAlignment = -4108

# Column widths
$ws.Columns.Item(2).ColumnWidth = 13.109375
$ws.Columns.Item(3).ColumnWidth = 15.44140625
$ws.Columns.Item(4).ColumnWidth = 18.21875
$ws.Columns.Item(5).ColumnWidth = 81.33203125

# Row height
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 13.8

# Selection
$ws.Range("E3").Select()
